$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Across the different user segments, a common pattern is that most individuals consume between two to three cups of coffee per day. Specifically, males and rural residents tend to drink around two cups daily, often timing their consumption with meals or social interactions. Urban dwellers and females typically consume three cups, with a structured routine of having coffee at set times throughout the day. The key difference lies in the slight increase in consumption among urban and female segments compared to males and rural residents. Overall, the user experience suggests a habitual and somewhat patterned approach to coffee consumption, with variations influenced by gender and living environment."

$ws.Range("B3").Value = "Across the user segments, a common pattern is the preference for making the first cup of coffee at home, as it is part of their morning routine. However, there is a tendency among urban and female users to buy coffee outside later in the day for convenience or as a treat. A key difference is that male users tend to buy coffee outside more frequently due to their busy schedules, while rural users predominantly make coffee at home, enjoying the brewing ritual. Overall, the user experience reflects a balance between the convenience of purchasing coffee outside and the comfort or ritual of making it at home."

$ws.Range("B4").Value = "Across all user segments, taste is consistently identified as the most important factor when choosing coffee. There is also a shared interest in quality, particularly from urban and female respondents who prefer to support local roasters. Price is a significant consideration for male and rural respondents, likely due to budget constraints such as being a student or living on a fixed income. Convenience is specifically noted by male respondents as an important factor. Overall, the user experience shows a strong emphasis on taste, with varying levels of concern for price and support for local businesses depending on the demographic."

$ws.Range("B5").Value = "The common pattern across the segments is that all interviewees believe coffee impacts their energy and productivity positively. The male, urban, and female segments emphasize that coffee significantly helps them stay focused and alert, particularly during demanding activities like exams or meetings. The rural segment, however, views coffee more as a comforting morning ritual rather than a necessity for boosting productivity. Key differences include the reliance on coffee; urban and female segments see it as essential for maintaining productivity, while the rural segment enjoys it more for leisure. Overall, coffee is perceived as a beneficial factor for energy and productivity, though its role varies from being a crucial tool for managing busy schedules to simply enhancing the morning experience."

$ws.Range("B6").Value = "The user research responses reveal that coffee consumption is intertwined with both personal and professional experiences across different segments. A common pattern is that coffee plays a significant role in memorable life events, whether they are personal or professional. For the male and urban segments, coffee is associated with career-related achievements, such as acing a test or securing a major client deal. In contrast, the rural segment associates coffee with personal life, specifically a cherished daily routine with a loved one. The female segment shares a similar professional anecdote as the urban segment, highlighting the importance of coffee in networking and career progression. Overall, coffee serves as a catalyst for important life moments, offering both personal enjoyment and professional opportunities."
